$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text / plain replacements (not numeric-looking, safe as-is) ---
$ws.Range("C2").Value  = "Hartmut"
$ws.Range("C3").Value  = "Mohaupt"

$ws.Range("D5").Value  = "KONTOSTAND AM 04.09.2024"

$ws.Range("B6").Value  = "07.09."
$ws.Range("C6").Value  = "08.09."
$ws.Range("D6").Value  = "PAYPAL LQLZYN"
$ws.Range("E6").Value  = "17,23-"

$ws.Range("B7").Value  = "08.09."
$ws.Range("C7").Value  = "09.09."
$ws.Range("D7").Value  = "PAYPAL PQYPMH"
$ws.Range("E7").Value  = "49,13-"

$ws.Range("B8").Value  = "12.09."
$ws.Range("C8").Value  = "13.09."
$ws.Range("D8").Value  = "RECHNUNG VODAFONE GMBH 37642085"
$ws.Range("E8").Value  = "37,69-"

$ws.Range("D12").Value = "KONTOSTAND AM 15.09.2024"
$ws.Range("E12").Value = "104,05-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 24.09.2024"

# --- B3 holds a purely numeric-looking card number. It must stay plain text
#     (inline/shared string) with its original cell style (s="8"), exactly as
#     before. Directly assigning a digit-only string turns it into a number,
#     and pre-formatting the cell as Text ("@") before assignment keeps it
#     textual but stamps the style with a "quote prefix" flag (a different
#     style id). To avoid both side effects: switch the cell to Text format,
#     set the value, then restore the original formatting by pasting the
#     format from an untouched neighbor cell that still carries style "8"
#     (e.g. D9, which this edit never touches).
$cell = $ws.Range("B3")
$styleDonor = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2570314725427075"
$styleDonor.Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
